$d = $word.ActiveDocument

function Replace-Literal($findText, $replaceText) {
    $range = $d.Content
    $found = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $findText"
    }
}

Replace-Literal "[court.name]" "`${court-name}"
Replace-Literal "[court.address]" "`${court-address}"
Replace-Literal "[plaintiff.name]" "`${plaintiff-name}"
Replace-Literal "[plaintiff.info]" "`${plaintiff-info}"
Replace-Literal "[plaintiff.representative]" "`${plaintiff-representative}"
Replace-Literal "[defendant.name]" "`${defendant-name}"
Replace-Literal "[defendant.info]" "`${defendant-info}"
Replace-Literal "[defendant.representative]" "`${defendant-representative}"
Replace-Literal "[counter]" "`${counter}"
Replace-Literal "[formelles]" "`${formelles}"
Replace-Literal "[materielles]" "`${materielles}"
Replace-Literal "[representative.name]" "`${representative-name}"

Write-Output "Done"
